$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before O (pushes old O..S -> P..T) and add the new
# "Year" ("سال") column header + per-row year values.
$ws.Columns("O").Insert()

$ws.Range("O1").Value = "سال"
$ws.Range("O3").Value = "۱۴۰۰"
$ws.Range("O4").Value = "۱۳۹۹"
$ws.Range("O6").Value = 2011
